# Update annotations for Sunsi Wu
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B48: was stored as text "1", should be a real number 1
$ws.Range("B48").Value = 1

# Append new row 49 with the new annotation record.
# Note: B49 must stay text "3" (not numeric), so prefix with an apostrophe
# to force Excel to keep it as a text entry instead of auto-converting.
$ws.Range("A49").Value = "Sunsi Wu"
$ws.Range("B49").Value = "'3"
$ws.Range("C49").Value = "无"
$ws.Range("D49").Value = "DFT"
$ws.Range("E49").Value = "EXP"
$ws.Range("F49").Value = "3bf3a8cd-f7a3-492e-815a-c1d9e74634b1"
$ws.Range("G49").Value = "ByCPHrgCW_annotated.xlsx"
$ws.Range("H49").Value = "The problem scenario states that the model/weights is private, but later on it ceases to be so (weights are not encrypted)."
